# Insert a new weekly price record at the top of the price-history block
# (row 559), pushing the existing rows 559-610 down to 560-611.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(559).Insert()

$ws.Range("A559").Value = 3
$ws.Range("B559").Value = "Femacal de La Calera"
$ws.Range("C559").Value = "Coquimbo"
$ws.Range("D559").Value = 45166
$ws.Range("E559").Value = 5
$ws.Range("F559").Value = 100112009
$ws.Range("G559").Value = "Acelga"
$ws.Range("H559").Value = "Sin especificar"
$ws.Range("I559").Value = "Primera"
$ws.Range("J559").Value = 210
$ws.Range("K559").Value = 3500
$ws.Range("L559").Value = 3800
$ws.Range("M559").Value = 3643
$ws.Range("N559").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O559").Value = "Provincia de Quillota"
$ws.Range("P559").Value = 607
$ws.Range("Q559").Value = 6
$ws.Range("R559").Value = "Hortaliza"
